$wb = $excel.ActiveWorkbook

# Sheet "2025" (sheet1) - row 2 updated investment-cost values
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 2900.628494009472
$ws.Range("E2").Value = 290490.7128553879
$ws.Range("G2").Value = 80959.2571266193
$ws.Range("I2").Value = 149402.1181152952
$ws.Range("L2").Value = 509988.6069102
$ws.Range("M2").Value = 112287.0813999
$ws.Range("N2").Value = 71616.34392528556
$ws.Range("O2").Value = 66890.96019342256

# Sheet "2030" (sheet2) - row 2 updated investment-cost values
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 2196.191807040655
$ws.Range("B2").Value = 35136.73892605113
$ws.Range("E2").Value = 164484.773501275
$ws.Range("I2").Value = 163495.084531733
$ws.Range("L2").Value = 90560.81452240903
$ws.Range("M2").Value = 61433.01601085002
$ws.Range("N2").Value = 19369.80367784133
$ws.Range("O2").Value = 11605.57876634009

# Sheet "2035" (sheet3) - row 2 updated investment-cost values
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 19330.2301468956
$ws.Range("B2").Value = 19544.39826649252
$ws.Range("E2").Value = 120007.9638621264
$ws.Range("I2").Value = 169968.1505399035
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 59687.40270289499
$ws.Range("N2").Value = 44280.58224493515
$ws.Range("O2").Value = 51872.05916375208
